$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 55.65000152587891
$ws.Range("E2").Value = 62.95000076293945
$ws.Range("F2").Value = 68.63999938964844
$ws.Range("G2").Value = 51.15999984741211
$ws.Range("H2").Value = 50476952
$ws.Range("I2").Value = "CYBR"

$ws.Range("D3").Value = 64
$ws.Range("E3").Value = 59.13000106811523
$ws.Range("F3").Value = 65.5
$ws.Range("G3").Value = 53.19200134277344
$ws.Range("H3").Value = 50476952
$ws.Range("I3").Value = "CYBR"

$ws.Range("D4").Value = 49.56000137329102
$ws.Range("E4").Value = 49.63999938964844
$ws.Range("F4").Value = 54.7400016784668
$ws.Range("G4").Value = 45.02999877929688
$ws.Range("H4").Value = 50476952
$ws.Range("I4").Value = "CYBR"

$ws.Range("D5").Value = 43.9900016784668
$ws.Range("E5").Value = 43.58000183105469
$ws.Range("F5").Value = 49.56000137329102
$ws.Range("G5").Value = 36.3380012512207
$ws.Range("H5").Value = 50476952
$ws.Range("I5").Value = "CYBR"

$ws.Range("D6").Value = 42.20999908447266
$ws.Range("E6").Value = 40.84000015258789
$ws.Range("F6").Value = 43.29999923706055
$ws.Range("G6").Value = 37
$ws.Range("H6").Value = 50476952
$ws.Range("I6").Value = "CYBR"

$ws.Range("D7").Value = 48.38000106811523
$ws.Range("E7").Value = 56.38999938964844
$ws.Range("F7").Value = 57.40999984741211
$ws.Range("G7").Value = 47.81999969482422
$ws.Range("H7").Value = 50476952
$ws.Range("I7").Value = "CYBR"

$ws.Range("D8").Value = 49.63999938964844
$ws.Range("E8").Value = 46.75
$ws.Range("F8").Value = 51.56000137329102
$ws.Range("G8").Value = 44.56999969482422
$ws.Range("H8").Value = 50476952
$ws.Range("I8").Value = "CYBR"

$ws.Range("D10").Value = 51.11999893188477
$ws.Range("E10").Value = 52.90999984741211
$ws.Range("F10").Value = 54.04000091552734
$ws.Range("G10").Value = 50.22000122070312
$ws.Range("H10").Value = 50476952
$ws.Range("I10").Value = "CYBR"

$ws.Range("D11").Value = 50.02000045776367
$ws.Range("E11").Value = 41.68000030517578
$ws.Range("F11").Value = 51.2599983215332
$ws.Range("G11").Value = 41.31999969482422
$ws.Range("H11").Value = 50476952
$ws.Range("I11").Value = "CYBR"

$ws.Range("D12").Value = 41.06000137329102
$ws.Range("E12").Value = 42.36999893188477
$ws.Range("F12").Value = 44.18999862670898
$ws.Range("G12").Value = 40.61999893188477
$ws.Range("H12").Value = 50476952
$ws.Range("I12").Value = "CYBR"

$ws.Range("D13").Value = 41.59999847412109
$ws.Range("E13").Value = 43.27999877929688
$ws.Range("F13").Value = 44.29999923706055
$ws.Range("G13").Value = 40.63000106811523
$ws.Range("H13").Value = 50476952
$ws.Range("I13").Value = "CYBR"

$ws.Range("D14").Value = 50.88000106811523
$ws.Range("E14").Value = 54.95999908447266
$ws.Range("F14").Value = 55.18999862670898
$ws.Range("G14").Value = 48.59999847412109
$ws.Range("H14").Value = 50476952
$ws.Range("I14").Value = "CYBR"

$ws.Range("D15").Value = 62.97000122070312
$ws.Range("E15").Value = 60.70999908447266
$ws.Range("F15").Value = 69.97000122070312
$ws.Range("G15").Value = 59.02000045776367
$ws.Range("H15").Value = 50476952
$ws.Range("I15").Value = "CYBR"

$ws.Range("D16").Value = 80.75
$ws.Range("E16").Value = 68.26000213623047
$ws.Range("F16").Value = 81.23999786376953
$ws.Range("G16").Value = 65.84999847412109
$ws.Range("H16").Value = 50476952
$ws.Range("I16").Value = "CYBR"

$ws.Range("D17").Value = 72.90000152587891
$ws.Range("E17").Value = 87.76000213623047
$ws.Range("F17").Value = 88.93699645996094
$ws.Range("G17").Value = 69.15000152587891
$ws.Range("H17").Value = 50476952
$ws.Range("I17").Value = "CYBR"

$ws.Range("D18").Value = 121
$ws.Range("E18").Value = 128.9299926757812
$ws.Range("F18").Value = 129.9299926757812
$ws.Range("G18").Value = 111.25
$ws.Range("H18").Value = 50476952
$ws.Range("I18").Value = "CYBR"

$ws.Range("D19").Value = 130.8099975585938
$ws.Range("E19").Value = 138.8800048828125
$ws.Range("F19").Value = 148.7400054931641
$ws.Range("G19").Value = 124.7699966430664
$ws.Range("H19").Value = 50476952
$ws.Range("I19").Value = "CYBR"

$ws.Range("D20").Value = 99.91000366210938
$ws.Range("E20").Value = 101.5800018310547
$ws.Range("F20").Value = 106.806999206543
$ws.Range("G20").Value = 94.3000030517578
$ws.Range("H20").Value = 50476952
$ws.Range("I20").Value = "CYBR"

$ws.Range("D21").Value = 118.1900024414062
$ws.Range("E21").Value = 138.2299957275391
$ws.Range("F21").Value = 142.9889984130859
$ws.Range("G21").Value = 117.120002746582
$ws.Range("H21").Value = 50476952
$ws.Range("I21").Value = "CYBR"

$ws.Range("D22").Value = 82.80999755859375
$ws.Range("E22").Value = 98.76000213623048
$ws.Range("F22").Value = 100.5400009155273
$ws.Range("G22").Value = 75.12999725341797
$ws.Range("H22").Value = 50476952
$ws.Range("I22").Value = "CYBR"

$ws.Range("D23").Value = 99.26000213623048
$ws.Range("E23").Value = 117.8399963378906
$ws.Range("F23").Value = 120.5
$ws.Range("G23").Value = 98.56999969482422
$ws.Range("H23").Value = 50476952
$ws.Range("I23").Value = "CYBR"

$ws.Range("D24").Value = 104.7799987792969
$ws.Range("E24").Value = 99.15000152587891
$ws.Range("F24").Value = 116.4349975585938
$ws.Range("G24").Value = 98.23999786376952
$ws.Range("H24").Value = 50476952
$ws.Range("I24").Value = "CYBR"

$ws.Range("D25").Value = 161.8600006103516
$ws.Range("E25").Value = 160.25
$ws.Range("F25").Value = 169.6999969482422
$ws.Range("G25").Value = 142.7799987792969
$ws.Range("H25").Value = 50476952
$ws.Range("I25").Value = "CYBR"

$ws.Range("D26").Value = 131.6999969482422
$ws.Range("E26").Value = 140.5
$ws.Range("F26").Value = 154.3699951171875
$ws.Range("G26").Value = 131.0500030517578
$ws.Range("H26").Value = 50476952
$ws.Range("I26").Value = "CYBR"

$ws.Range("D27").Value = 130.4799957275391
$ws.Range("E27").Value = 142.0299987792969
$ws.Range("F27").Value = 147.3500061035156
$ws.Range("G27").Value = 129.8999938964844
$ws.Range("H27").Value = 50476952
$ws.Range("I27").Value = "CYBR"

$ws.Range("D28").Value = 158.5200042724609
$ws.Range("E28").Value = 180.1100006103516
$ws.Range("F28").Value = 187.8699951171875
$ws.Range("G28").Value = 152.9299926757812
$ws.Range("H28").Value = 50476952
$ws.Range("I28").Value = "CYBR"

$ws.Range("D29").Value = 173.7400054931641
$ws.Range("E29").Value = 137.1499938964844
$ws.Range("F29").Value = 173.7400054931641
$ws.Range("G29").Value = 125.0199966430664
$ws.Range("H29").Value = 50476952
$ws.Range("I29").Value = "CYBR"

$ws.Range("D30").Value = 169.7400054931641
$ws.Range("E30").Value = 157.1399993896484
$ws.Range("F30").Value = 180.0099945068359
$ws.Range("G30").Value = 155.0899963378906
$ws.Range("H30").Value = 50476952
$ws.Range("I30").Value = "CYBR"

$ws.Range("D31").Value = 127.5100021362305
$ws.Range("E31").Value = 130.1300048828125
$ws.Range("F31").Value = 144.8200073242188
$ws.Range("G31").Value = 123.1650009155273
$ws.Range("H31").Value = 50476952
$ws.Range("I31").Value = "CYBR"

$ws.Range("D32").Value = 150.5599975585938
$ws.Range("E32").Value = 156.9100036621094
$ws.Range("F32").Value = 162.6100006103516
$ws.Range("G32").Value = 132.5800018310547
$ws.Range("H32").Value = 50476952
$ws.Range("I32").Value = "CYBR"

$ws.Range("D33").Value = 131.6999969482422
$ws.Range("E33").Value = 140.8800048828125
$ws.Range("F33").Value = 141.6399993896484
$ws.Range("G33").Value = 113.1900024414062
$ws.Range("H33").Value = 50476952
$ws.Range("I33").Value = "CYBR"

$ws.Range("D34").Value = 146.9299926757812
$ws.Range("E34").Value = 124.5999984741211
$ws.Range("F34").Value = 146.9299926757812
$ws.Range("G34").Value = 121.4400024414062
$ws.Range("H34").Value = 50476952
$ws.Range("I34").Value = "CYBR"

$ws.Range("D35").Value = 156.0099945068359
$ws.Range("E35").Value = 166.0099945068359
$ws.Range("F35").Value = 169.3399963378906
$ws.Range("G35").Value = 152.3800048828125
$ws.Range("H35").Value = 50476952
$ws.Range("I35").Value = "CYBR"

$ws.Range("D36").Value = 164.8800048828125
$ws.Range("E36").Value = 163.6399993896484
$ws.Range("F36").Value = 171.6199951171875
$ws.Range("G36").Value = 152.0299987792969
$ws.Range("H36").Value = 50476952
$ws.Range("I36").Value = "CYBR"

$ws.Range("D37").Value = 214.8399963378907
$ws.Range("E37").Value = 233.479995727539
$ws.Range("F37").Value = 241.3619995117188
$ws.Range("G37").Value = 205.5899963378907
$ws.Range("H37").Value = 50476952
$ws.Range("I37").Value = "CYBR"

$ws.Range("D38").Value = 264.989990234375
$ws.Range("E38").Value = 239.25
$ws.Range("F38").Value = 267.5199890136719
$ws.Range("G38").Value = 230.3800048828125
$ws.Range("H38").Value = 50476952
$ws.Range("I38").Value = "CYBR"

$ws.Range("D39").Value = 272.8999938964844
$ws.Range("E39").Value = 256.3800048828125
$ws.Range("F39").Value = 284.1000061035156
$ws.Range("G39").Value = 249.0500030517578
$ws.Range("H39").Value = 50476952
$ws.Range("I39").Value = "CYBR"

$ws.Range("D40").Value = 294.5499877929688
$ws.Range("E40").Value = 276.5199890136719
$ws.Range("F40").Value = 308.6300048828125
$ws.Range("G40").Value = 269
$ws.Range("H40").Value = 50476952
$ws.Range("I40").Value = "CYBR"

$ws.Range("D41").Value = 335.6099853515625
$ws.Range("E41").Value = 370.9800109863281
$ws.Range("F41").Value = 378.6199951171875
$ws.Range("G41").Value = 331.6400146484375
$ws.Range("H41").Value = 50476952
$ws.Range("I41").Value = "CYBR"

$ws.Range("D42").Value = 338.1900024414062
$ws.Range("E42").Value = 352.1600036621094
$ws.Range("F42").Value = 359.5350036621094
$ws.Range("G42").Value = 288.6300048828125
$ws.Range("H42").Value = 50476952
$ws.Range("I42").Value = "CYBR"

$ws.Range("D43").Value = 406.25
$ws.Range("E43").Value = 411.4700012207031
$ws.Range("F43").Value = 452
$ws.Range("G43").Value = 370.1600036621094
$ws.Range("H43").Value = 50476952
$ws.Range("I43").Value = "CYBR"
